$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 ---
$ws.Range("C9").Value = 0.36087121474350697
$ws.Range("D9").Value = 0.56937035610904596
$ws.Range("G9").Value = 235
$ws.Range("I9").Value = 0.49907721187931903
$ws.Range("J9").Value = 0.61139622805221505
$ws.Range("M9").Value = 0.0491249163655313
$ws.Range("N9").Value = 0.79065971708551697
$ws.Range("Q9").Value = 278
$ws.Range("S9").Value = 0.0574302079812866
$ws.Range("T9").Value = 0.85263306574424202

# --- Row 24 ---
$ws.Range("C24").Value = 0.205770491432575
$ws.Range("D24").Value = 0.693555671227017
$ws.Range("G24").Value = 900
$ws.Range("I24").Value = 0.35461114690213802
$ws.Range("J24").Value = 0.75805200671086503
$ws.Range("M24").Value = 0.0341274160573958
$ws.Range("N24").Value = 0.399974400676194
$ws.Range("Q24").Value = 1191
$ws.Range("S24").Value = 0.0392279030920024
$ws.Range("T24").Value = 0.41294973356557202

# --- Row 25 ---
$ws.Range("C25").Value = 0.32603692591421801
$ws.Range("D25").Value = 0.540485768850547
$ws.Range("E25").Value = 0.54928383533993297
$ws.Range("F25").Value = 0.41310262328491898
$ws.Range("G25").Value = 1061
$ws.Range("I25").Value = 0.476610058523047
$ws.Range("J25").Value = 0.59136620874102797
$ws.Range("K25").Value = 0.802958745157622
$ws.Range("L25").Value = 0.60388517315260104
$ws.Range("M25").Value = 0.0287785686575152
$ws.Range("N25").Value = 0.53279319371673095
$ws.Range("O25").Value = 0.45623229977916302
$ws.Range("P25").Value = 0.0720734073534815
$ws.Range("Q25").Value = 1115
$ws.Range("S25").Value = 0.0353344040288236
$ws.Range("T25").Value = 0.55892886314922097
$ws.Range("U25").Value = 0.56016324519970795
$ws.Range("V25").Value = 0.0884919234680863

# --- Selection moves from J27 to M24 ---
$ws.Range("M24").Select()
